$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 608; existing rows 608-619 shift down to 610-621.
$ws.Range("A608:A609").EntireRow.Insert()

# Populate the two newly inserted rows (608-609) with the new price records.
# Columns A,B,C,E,F,G,H,I,J,K,T are identical to the surrounding rows for this
# market/product (Terminal Hortofrutícola Agro Chillán - Frutilla).
$ws.Cells.Item(608, 1).Value = 7
$ws.Cells.Item(608, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(608, 3).Value = "Ñuble"
$ws.Cells.Item(608, 4).Value = 45239
$ws.Cells.Item(608, 5).Value = 16
$ws.Cells.Item(608, 6).Value = "Fruta"
$ws.Cells.Item(608, 7).Value = 100101
$ws.Cells.Item(608, 8).Value = "Berries"
$ws.Cells.Item(608, 9).Value = 100112025
$ws.Cells.Item(608, 10).Value = "Frutilla"
$ws.Cells.Item(608, 11).Value = "Sin especificar"
$ws.Cells.Item(608, 12).Value = "Especial"
$ws.Cells.Item(608, 13).Value = 100
$ws.Cells.Item(608, 14).Value = 15000
$ws.Cells.Item(608, 15).Value = 15000
$ws.Cells.Item(608, 16).Value = 15000
$ws.Cells.Item(608, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(608, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(608, 19).Value = 2143
$ws.Cells.Item(608, 20).Value = 7

$ws.Cells.Item(609, 1).Value = 7
$ws.Cells.Item(609, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(609, 3).Value = "Ñuble"
$ws.Cells.Item(609, 4).Value = 45239
$ws.Cells.Item(609, 5).Value = 16
$ws.Cells.Item(609, 6).Value = "Fruta"
$ws.Cells.Item(609, 7).Value = 100101
$ws.Cells.Item(609, 8).Value = "Berries"
$ws.Cells.Item(609, 9).Value = 100112025
$ws.Cells.Item(609, 10).Value = "Frutilla"
$ws.Cells.Item(609, 11).Value = "Sin especificar"
$ws.Cells.Item(609, 12).Value = "Primera"
$ws.Cells.Item(609, 13).Value = 100
$ws.Cells.Item(609, 14).Value = 13000
$ws.Cells.Item(609, 15).Value = 13000
$ws.Cells.Item(609, 16).Value = 13000
$ws.Cells.Item(609, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(609, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(609, 19).Value = 1857
$ws.Cells.Item(609, 20).Value = 7

# Ensure the date cells keep the same date number format used elsewhere in column D.
$ws.Range("D608:D609").NumberFormat = $ws.Range("D610").NumberFormat
